$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '75.920.80'
$ws.Range("E2").Value = '  +9.79%  '
$ws.Range("D3").Value = '2.694.02'
$ws.Range("E3").Value = '  +11.53%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '189.20'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '591.75'
$ws.Range("E6").Value = '  +5.36%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  +5.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.199'
$ws.Range("E9").Value = '  +18.63%  '
$ws.Range("D10").Value = '2.690.86'
$ws.Range("E10").Value = '  +11.54%  '
$ws.Range("E12").Value = '  +8.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.75'
$ws.Range("E13").Value = '  +2.35%  '
$ws.Range("D14").Value = '75.722.34'
$ws.Range("E14").Value = '  +9.67%  '
$ws.Range("D15").Value = '3.189.25'
$ws.Range("E15").Value = '  +11.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000191'
$ws.Range("E16").Value = '  +8.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.72'
$ws.Range("E17").Value = '  +12.04%  '
$ws.Range("D18").Value = '2.686.62'
$ws.Range("E18").Value = '  +10.77%  '
$ws.Range("E19").Value = '  +33.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.13'
$ws.Range("E20").Value = '  +12.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.28'
$ws.Range("E21").Value = '  +10.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.29'
$ws.Range("E22").Value = '  +17.19%  '
$ws.Range("E23").Value = '  +5.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.28'
$ws.Range("E24").Value = '  +4.79%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.76'
$ws.Range("E26").Value = '  +7.45%  '
$ws.Range("E27").Value = '  +10.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.53'
$ws.Range("E28").Value = '  +13.00%  '
$ws.Range("D29").Value = '2.830.11'
$ws.Range("E29").Value = '  +11.35%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").Value = '0.0₃0968'
$ws.Range("E31").Value = '  +14.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '522.88'
$ws.Range("E32").Value = '  +16.03%  '
$ws.Range("E33").Value = '  +14.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.89'
$ws.Range("E34").Value = '  +7.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.78'
$ws.Range("E35").Value = '  +10.60%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  +9.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.04'
$ws.Range("E38").Value = '  +2.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.38'
$ws.Range("E39").Value = '  +6.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.40'
$ws.Range("E40").Value = '  +1.69%  '
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.07'
$ws.Range("E42").Value = '  +15.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '171.92'
$ws.Range("E43").Value = '  +27.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.72'
$ws.Range("E44").Value = '  +13.25%  '
$ws.Range("E45").Value = '  +10.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.21'
$ws.Range("E46").Value = '  +11.69%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.41'
$ws.Range("E47").Value = '  +15.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '39.40'
$ws.Range("E48").Value = '  +4.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0854'
$ws.Range("E49").Value = '  +18.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.68'
$ws.Range("E50").Value = '  +8.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.544'
$ws.Range("E51").Value = '  +11.71%  '
